# Add 2022-Q4 data
#
# 1) "总计" (summary) sheet: insert a new top data row for 2022-Q4 and push
#    the existing 2022-Q3 / 2022-Q1 rows down by one.
# 2) Add a brand-new "2022-Q4" worksheet (positioned right after "总计",
#    before "2022-Q3") holding the per-fund breakdown for the new quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) 总计 sheet — shift rows 2,3 down to 3,4 and write the new row 2
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Make room / carry formatting: copy row 3 (with its style) down onto the
# new row 4, then row 2 down onto row 3.
$summary.Range("A3:D3").Copy($summary.Range("A4:D4"))
$summary.Range("A2:D2").Copy($summary.Range("A3:D3"))

# Row 4 <- old 2022-Q1 row (values already copied above; just fix index)
$summary.Cells.Item(4, 1).Value = 2

# Row 3 <- old 2022-Q3 row (values already copied above; just fix index)
$summary.Cells.Item(3, 1).Value = 1

# Row 2 <- brand new 2022-Q4 row
$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q4"
$summary.Cells.Item(2, 3).Value = 3
$summary.Cells.Item(2, 4).Value = 0.15

# ---------------------------------------------------------------------
# 2) New "2022-Q4" worksheet with the per-fund holdings
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")

# Duplicate the "2022-Q3" sheet (same column layout/styling) and drop the
# copy in right before it; rename to "2022-Q4".
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# The template only has one data row (row 2) — extend formatting down to
# rows 3 and 4 for the two extra funds.
$q4.Range("A2:H2").Copy($q4.Range("A3:H3"))
$q4.Range("A2:H2").Copy($q4.Range("A4:H4"))

# Row 2: 009623 长城创新驱动混合A
$q4.Cells.Item(2, 1).Value = 0
$q4.Cells.Item(2, 2).Value = "'009623"
$q4.Cells.Item(2, 3).Value = "长城创新驱动混合A"
$q4.Cells.Item(2, 4).Value = "'4.20"
$q4.Cells.Item(2, 5).Value = "'91.95"
$q4.Cells.Item(2, 6).Value = "'2.93"
$q4.Cells.Item(2, 7).Value = "'0.1231"
$q4.Cells.Item(2, 8).Value = 10

# Row 3: 519097 新华中小市值优选混合
$q4.Cells.Item(3, 1).Value = 1
$q4.Cells.Item(3, 2).Value = "'519097"
$q4.Cells.Item(3, 3).Value = "新华中小市值优选混合"
$q4.Cells.Item(3, 4).Value = "'0.66"
$q4.Cells.Item(3, 5).Value = "'70.51"
$q4.Cells.Item(3, 6).Value = "'3.59"
$q4.Cells.Item(3, 7).Value = "'0.0237"
$q4.Cells.Item(3, 8).Value = 6

# Row 4: 017458 长城创新驱动混合C
$q4.Cells.Item(4, 1).Value = 2
$q4.Cells.Item(4, 2).Value = "'017458"
$q4.Cells.Item(4, 3).Value = "长城创新驱动混合C"
$q4.Cells.Item(4, 4).Value = "'0.00"
$q4.Cells.Item(4, 5).Value = "'91.95"
$q4.Cells.Item(4, 6).Value = "'2.93"
$q4.Cells.Item(4, 7).Value = 0
$q4.Cells.Item(4, 8).Value = 10

Write-Output "2022-Q4 sheet added"
